$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Three new arrival rows appended after the existing 72 data rows (rows 73-75),
# continuing the running "NUMBER" sequence (72, 73, 74) for Saturday, Jan 14.

$newRows = @(
    @{ Num = 72; Date = "Saturday, Jan 14"; Time = "9:15 PM";  Flight = "W95153"; From = "London";           Short = "(LTN)"; Airline = "Wizz Air ";                                  Model = "A320"; AircraftId = "(G-WUKF)"; Status = "8:48 PM";  Diff = "0 hours, -27 minutes" },
    @{ Num = 73; Date = "Saturday, Jan 14"; Time = "9:40 PM";  Flight = "FR3472";  From = "London";           Short = "(LTN)"; Airline = "Ryanair ";                                    Model = "B738"; AircraftId = "(EI-EFJ)"; Status = "9:54 PM";  Diff = "0 hours, 14 minutes" },
    @{ Num = 74; Date = "Saturday, Jan 14"; Time = "11:05 PM"; Flight = "PQ7551";  From = "Sharm el-Sheikh"; Short = "(SSH)"; Airline = "SkyUp Airlines (The Power Of Freedom Livery) "; Model = "B738"; AircraftId = "(UR-SQM)"; Status = "11:11 PM"; Diff = "0 hours, 6 minutes" }
)

$r = 73
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row.Num
    $ws.Cells.Item($r, 2).Value = $row.Date
    $ws.Cells.Item($r, 3).Value = $row.Time
    $ws.Cells.Item($r, 4).Value = $row.Flight
    $ws.Cells.Item($r, 5).Value = $row.From
    $ws.Cells.Item($r, 6).Value = $row.Short
    $ws.Cells.Item($r, 7).Value = $row.Airline
    $ws.Cells.Item($r, 8).Value = $row.Model
    $ws.Cells.Item($r, 9).Value = $row.AircraftId
    $ws.Cells.Item($r, 10).Value = $row.Status
    # Column K (11) is always blank in this sheet (a pre-existing layout quirk
    # where the header's "DIFFERENCE" label sits in K1 but the data rows hold
    # the value one column over, in L) - materialize the blank cell to match
    # every other data row without introducing a new cell style.
    $ws.Cells.Item($r, 11).Font.Bold = $false
    $ws.Cells.Item($r, 12).Value = $row.Diff
    $ws.Cells.Item($r, 13).Font.Bold = $false
    $r = $r + 1
}
